$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row updates derived from the cryptos price-refresh diff.
# Columns: B=Coin name, C=Link, D=Price, E=Volume(1h).
# D-column values that are plain numeric strings need NumberFormat forced
# to Text before the write (then restored to General) so Excel keeps them
# as text -- matching the "71.214.04"-style prices elsewhere in the sheet --
# instead of silently converting them to a number.

$ws.Range("D2").Value = "71.213.03"
$ws.Range("E2").Value = "  +0.40%  "

$ws.Range("D3").Value = "3.814.79"
$ws.Range("E3").Value = "  -0.87%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "705.45"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +1.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.95"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -0.06%  "

$ws.Range("D7").Value = "3.813.73"
$ws.Range("E7").Value = "  -0.86%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("E9").Value = "  +0.24%  "

$ws.Range("E10").Value = "  -0.11%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.66"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +5.48%  "

$ws.Range("E12").Value = "  +0.69%  "

$ws.Range("E13").Value = "  -1.43%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.04"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -0.35%  "

$ws.Range("D15").Value = "4.458.79"
$ws.Range("E15").Value = "  -0.87%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "71.170.44"
$ws.Range("E16").Value = "  +0.27%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.742.15"
$ws.Range("E17").Value = "  -2.73%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.56"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +1.12%  "

$ws.Range("E19").Value = "  -0.01%  "

$ws.Range("E20").Value = "  -0.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "517.50"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +4.57%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.68"
$ws.Range("D22").NumberFormat = "General"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.723"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +0.71%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.30"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -0.46%  "

$ws.Range("E25").Value = "  -1.92%  "

$ws.Range("D26").Value = "3.966.10"
$ws.Range("E26").Value = "  -0.91%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.04"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -1.33%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.41"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -1.49%  "

$ws.Range("E29").Value = "  +0.24%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.04"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -3.31%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.04"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -2.79%  "

$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.40"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -1.35%  "

$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.25"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -0.43%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.18"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -0.98%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.172"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -4.15%  "

$ws.Range("E36").Value = "  +0.08%  "

$ws.Range("D37").Value = "3.777.10"
$ws.Range("E37").Value = "  -0.77%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -0.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.39"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -0.14%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.98"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -0.23%  "

$ws.Range("E42").Value = "  -1.78%  "

$ws.Range("E43").Value = "  -1.59%  "

$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "171.14"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +4.52%  "

$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +0.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.000313"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +1.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "49.55"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +1.67%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "422.23"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +3.47%  "

$ws.Range("E50").Value = "  +0.38%  "

$ws.Range("E51").Value = "  -1.46%  "
